$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells keep their exact string representation
# (avoids Excel auto-converting numeric-looking strings to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.687.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.505.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.92"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.503.61"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000217"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.094.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.502.55"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.561.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.56"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.87"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.67"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.629"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.16"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.642.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000125"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.38%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.62"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.15"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.494.85"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.97"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "178.74"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0892"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.45"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.84%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.30"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.990"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.61%  "
